# Loan RBI, Variable Instalments
#
# - On "Transactions", move the selection to D22 (it is no longer the
#   active/selected tab).
# - On "Repayment schedule", insert a new (blank) column before the old
#   column N ("Late"), which pushes "Late", the blank "heading" column and
#   "Outstanding" one column to the right. Re-apply the width that the new
#   column should carry (inherited from its left neighbour, "In Advance").
# - Make "Repayment schedule" the active sheet/tab and move its selection
#   to L15.

$wb = $excel.ActiveWorkbook

$wsTransactions = $wb.Worksheets.Item("Transactions")
[void]$wsTransactions.Range("D22").Select()

$wsRepayment = $wb.Worksheets.Item("Repayment schedule")
$wsRepayment.Columns("N").Insert()
$wsRepayment.Columns("N").ColumnWidth = $wsRepayment.Columns("M").ColumnWidth

[void]$wsRepayment.Activate()
[void]$wsRepayment.Range("L15").Select()
